$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Weekly Quantity" ---
$ws1 = $wb.Worksheets.Item("Weekly Quantity")

# Update B39: 18 -> 16
$ws1.Range("B39").Value = 16

# Delete entire row 40 (45368.99999999999 / 4), shifting subsequent rows up
$ws1.Rows.Item(40).Delete()

# --- Sheet 2: "Monthly Trend" ---
$ws2 = $wb.Worksheets.Item("Monthly Trend")

# Update B14: 22 -> 16
$ws2.Range("B14").Value = 16
